$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new columns I (I0) and J (IF) with the header formatting copied from column H
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 / IF values for data rows 2-51 (row, I0, IF)
$data = @(
    @(2, 7, 7),
    @(3, 4, 4),
    @(4, 3, 3),
    @(5, 8, 8),
    @(6, 1, 3),
    @(7, 7, 8),
    @(8, 6, 8),
    @(9, 10, 10),
    @(10, 7, 7),
    @(11, 8, 8),
    @(12, 7, 7),
    @(13, 6, 7),
    @(14, 9, 9),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 6, 7),
    @(19, 7, 8),
    @(20, 8, 8),
    @(21, 9, 9),
    @(22, 8, 9),
    @(23, 6, 6),
    @(24, 8, 8),
    @(25, 6, 6),
    @(26, 12, 12),
    @(27, 5, 5),
    @(28, 6, 6),
    @(29, 7, 7),
    @(30, 6, 6),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 6, 7),
    @(34, 8, 8),
    @(35, 6, 6),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 6, 6),
    @(39, 6, 7),
    @(40, 7, 7),
    @(41, 6, 7),
    @(42, 7, 7),
    @(43, 9, 9),
    @(44, 9, 9),
    @(45, 6, 6),
    @(46, 8, 8),
    @(47, 7, 8),
    @(48, 8, 8),
    @(49, 7, 8),
    @(50, 8, 8),
    @(51, 8, 8)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
